$d = $word.ActiveDocument

# 1. Replace the text of the first "mushroom" bullet with the new wording.
$d.Content.Find.Execute(
    "Make sure mushroom does not disappear when coming out of block.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Make sure that mushroom only collides with ground layers. Should go through creatures or player (of course as soon as player touches it, it should disappear and Mario should grow).",
    2
)

# 2. Remove the whole paragraph "Make sure that mushroom doesn't get stuck in cracks or bumps in the ground."
#    (including its trailing paragraph mark so the paragraph itself disappears).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Make sure that mushroom*get stuck in cracks or bumps*") {
        $p.Range.Delete()
        break
    }
}
